# Insert a new weekly price record at row 23 ("Fruta / hortaliza, semanal"
# commit). All existing rows from 23 down to 104 shift down by one (to
# 24..105), which is how the source diff shows every subsequent row's
# values "moving" to the next row while the sheet's used range grows from
# A1:R104 to A1:R105.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push row 23 (and everything below it) down by one row.
$ws.Rows("23:23").Insert()

# Populate the newly opened row 23 with the new record.
$ws.Cells.Item(23, 1).Value = 11
$ws.Cells.Item(23, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(23, 3).Value = "Bíobío"
$ws.Cells.Item(23, 4).Value = "2022-02-16"
$ws.Cells.Item(23, 5).Value = 8
$ws.Cells.Item(23, 6).Value = 100112043
$ws.Cells.Item(23, 7).Value = "Pepino ensalada"
$ws.Cells.Item(23, 8).Value = "Sin especificar"
$ws.Cells.Item(23, 9).Value = "Primera"
$ws.Cells.Item(23, 10).Value = 100
$ws.Cells.Item(23, 11).Value = 14000
$ws.Cells.Item(23, 12).Value = 15000
$ws.Cells.Item(23, 13).Value = 14500
$ws.Cells.Item(23, 14).Value = "$/caja 60 unidades"
$ws.Cells.Item(23, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(23, 16).Value = 242
$ws.Cells.Item(23, 17).Value = 60
$ws.Cells.Item(23, 18).Value = "Hortaliza"
